$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly (new text is longer)
$ws.Columns.Item(3).ColumnWidth = 56.140625

# Add new row 4: TwoSum exercise
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "TwoSum"
$ws.Range("C4").Value = "Given an array, return the indicies who's value sums to a target"
$ws.Range("D4").Value = "Nested for loop"
$ws.Range("E4").Value = 10

# Add row 5 with just the next index number
$ws.Range("A5").Value = 4

# Match the selection left after editing
$ws.Range("B5").Select()
